# The commit swaps the two theme parts in the package: the theme that was
# "Office Theme" (theme1.xml) and the theme that was "Integral" (theme2.xml)
# trade places, so the deck's live design (driven by the slide master's
# theme) ends up using the old "Office Theme" color values while the
# "Integral" palette moves to the part that is only referenced by the notes
# master.
#
# Through the PowerPoint object model, the live/applied theme is reached via
# SlideMaster.Theme (Slide.Design / NotesMaster.Theme / etc. all resolve back
# to this same applied theme). We update its color scheme, entry by entry, to
# the target ("Office Theme") RGB values so the applied design matches the
# post-commit palette.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colorScheme = $theme.ThemeColorScheme

# PowerPoint COM RGB() packs colors as r + g*256 + b*65536 (BGR word order).
# Target values below are the "Office Theme" palette that theme1.xml held
# before the edit (and that the applied theme must hold afterwards):
#   dk1=000000 lt1=FFFFFF dk2=44546A lt2=E7E6E6 accent1=5B9BD5
#   accent2=ED7D31 accent3=A5A5A5 accent4=FFC000 accent5=4472C4
#   accent6=70AD47 hlink=0563C1 folHlink=954F72
$colorScheme.Item(1).RGB = 0          # dk1      -> 000000
$colorScheme.Item(2).RGB = 16777215   # lt1      -> FFFFFF
$colorScheme.Item(3).RGB = 6968388    # dk2      -> 44546A
$colorScheme.Item(4).RGB = 15132391   # lt2      -> E7E6E6
$colorScheme.Item(5).RGB = 13998939   # accent1  -> 5B9BD5
$colorScheme.Item(6).RGB = 3243501    # accent2  -> ED7D31
$colorScheme.Item(7).RGB = 10855845   # accent3  -> A5A5A5
$colorScheme.Item(8).RGB = 49407      # accent4  -> FFC000
$colorScheme.Item(9).RGB = 12874308   # accent5  -> 4472C4
$colorScheme.Item(10).RGB = 4697456   # accent6  -> 70AD47
$colorScheme.Item(11).RGB = 12673797  # hlink    -> 0563C1
$colorScheme.Item(12).RGB = 7491477   # folHlink -> 954F72
